# Splits the two inline "{m" / "{m:" field-opening runs into a separate
# "{" run plus the remaining text, matching the
# TokenIteratorFieldRewriterSplit migration (commit message).
#
# Before: <w:r><w:t>{m</w:t></w:r>        (run immediately before ":v.name}")
#         <w:r><w:t>{m:</w:t></w:r>       (run immediately before "endfor}")
# After:  <w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r>
#         <w:r><w:t>{</w:t></w:r><w:r><w:t>m:</w:t></w:r>
#
# We locate each split point with Find (anchored on unique surrounding
# text) and force a run boundary right after the leading "{" by inserting
# a paragraph mark there and immediately deleting it again -- this leaves
# the two halves of the text as separate runs without adding any stray
# run-formatting (w:rPr) to either half.

$d = $word.ActiveDocument

function Split-RunAfterBrace([string]$anchorText) {
    $found = $d.Content
    $ok = $found.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $anchorText"
        return
    }
    $splitPos = $found.Start + 1
    $splitRange = $d.Range($splitPos, $splitPos)
    $splitRange.InsertParagraphAfter()
    $markRange = $d.Range($splitPos, $splitPos + 1)
    $markRange.Delete()
}

# First occurrence to split: the "{m" run directly before ":v.name}".
Split-RunAfterBrace("{m:v.name}")

# Second occurrence to split: the "{m:" run directly before "endfor}".
Split-RunAfterBrace("{m:endfor}")

Write-Output "Done splitting field-opening runs."
